# Update column G ("K") values on Sheet1 with freshly calculated strikeout
# counts (s_vals) now that the column represents K instead of Strike#.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 4
    4  = 9
    5  = 7
    6  = 7
    7  = 8
    8  = 5
    9  = 4
    10 = 5
    11 = 4
    12 = 4
    13 = 3
    14 = 6
    15 = 6
    16 = 5
    17 = 7
    18 = 4
    19 = 3
    20 = 5
    21 = 7
    22 = 7
    23 = 5
    24 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
